$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating point update on the execution timestamp of row 8
$ws.Range("A8").Value = 45875.37517995371

# Append the new reading row (row 9)
$ws.Range("A9").Value = 45875.41708634486
$ws.Range("B9").Value = 2025
$ws.Range("C9").Value = 23
$ws.Range("D9").Value = 18.06
$ws.Range("E9").Value = 83.13
$ws.Range("F9").Value = 457.03
$ws.Range("G9").Value = 7.54
$ws.Range("H9").Value = "ESE"
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "10:00:36"

# Match the date/time style used by the other timestamp cells in column A
$ws.Range("A9").NumberFormat = $ws.Range("A8").NumberFormat
